# Auto-generated from the authoritative OOXML diff.
# Updates the cryptos price/volume table cells that changed between commits.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.920.68"
$ws.Range("E2").Value = "  +4.26%  "
$ws.Range("D3").Value = "1.914.08"
$ws.Range("E3").Value = "  +1.50%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "248.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.11%  "
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "48.03"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +10.32%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.375"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "58.22"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.53%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0760"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.08%  "
$ws.Range("E12").Value = "  +2.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +12.83%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.827"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +7.36%  "
$ws.Range("D15").Value = "2.193.07"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.59%  "
$ws.Range("D17").Value = "1.900.56"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "36.943.10"
$ws.Range("E18").Value = "  +4.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "74.47"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.42%  "
$ws.Range("D20").Value = "0.0₃0857"
$ws.Range("E20").Value = "  +3.56%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.65"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "251.40"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.16"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.01%  "
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.45"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -7.22%  "
$ws.Range("E26").Value = "  +2.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "167.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.37%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "18.69"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.128"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.55%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.57"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0610"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.28%  "
$ws.Range("B33").Value = "Kaspa"
$ws.Range("C33").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0904"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +25.68%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.29"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.66%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.91"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.29%  "
$ws.Range("E36").Value = "  +0.07%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.32"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +52.66%  "
$ws.Range("E38").Value = "  +1.72%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.881"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.88%  "
$ws.Range("E40").Value = "  +1.80%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "104.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.53%  "
$ws.Range("E42").Value = "  +3.86%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "17.51"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.78%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.87"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +19.53%  "
$ws.Range("E45").Value = "  +2.13%  "
$ws.Range("D46").Value = "1.356.34"
$ws.Range("E46").Value = "  +3.53%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.41"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.73%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0829"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.37%  "
$ws.Range("E49").Value = "  +2.93%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.40"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "3.77"
$ws.Range("D51").Style = "Normal"
